$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K values (column G) per row, recomputed from regenerated save_data
$kValues = @(
    0, 0, 1, 1, 3, 1, 0, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 2, 1, 2, 1, 0, 1, 1, 2, 1, 0, 3, 0, 1, 1, 0, 1, 1, 0, 1, 1, 0, 0, 2, 1, 0, 0, 0, 2, 0, 0, 3, 0, 0, 0, 1, 1, 2, 0, 1, 2, 1, 0, 1, 2, 3, 3, 1, 2, 0, 1, 2, 3, 0, 1, 2, 1, 0, 3, 0, 1, 2, 2, 2
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

